# Apply cryptos list update (Fri Jan 19 17:33:43 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DirectValue($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-TextValue($cellRef, $val) {
    # Force the value to be stored as text even though it looks numeric
    # (e.g. "309.35"), then strip the temporary text format so the cell
    # keeps its original (default) style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-DirectValue "D2" "41.473.18"
Set-DirectValue "E2" "  -1.59%  "
Set-DirectValue "D3" "2.440.61"
Set-DirectValue "E3" "  -2.08%  "
Set-DirectValue "E4" "  +1.15%  "
Set-TextValue "D5" "309.35"
Set-DirectValue "E5" "  -1.15%  "
Set-TextValue "D6" "89.44"
Set-DirectValue "E6" "  -6.73%  "
Set-TextValue "D7" "0.531"
Set-DirectValue "E7" "  -4.58%  "
Set-DirectValue "E8" "  +0.95%  "
Set-TextValue "D9" "0.480"
Set-DirectValue "E9" "  -6.57%  "
Set-TextValue "D10" "31.56"
Set-DirectValue "E10" "  -8.26%  "
Set-DirectValue "E11" "  -3.18%  "
Set-DirectValue "E12" "  -0.22%  "
Set-DirectValue "D13" "2.816.21"
Set-DirectValue "E13" "  -1.99%  "
Set-TextValue "D14" "6.65"
Set-DirectValue "E14" "  -5.95%  "
Set-DirectValue "B15" "Chainlink"
Set-DirectValue "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "14.97"
Set-DirectValue "E15" "  +0.31%  "
Set-DirectValue "B16" "WrappedEther"
Set-DirectValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-DirectValue "D16" "2.406.21"
Set-DirectValue "E16" "  -2.68%  "
Set-DirectValue "E17" "  -5.88%  "
Set-DirectValue "D18" "41.136.41"
Set-DirectValue "E18" "  -2.27%  "
Set-TextValue "D19" "6.12"
Set-DirectValue "E19" "  -5.37%  "
Set-DirectValue "D20" "0.0₃0897"
Set-DirectValue "E20" "  -3.16%  "
Set-TextValue "D21" "68.77"
Set-DirectValue "E21" "  -0.94%  "
Set-TextValue "D22" "10.66"
Set-DirectValue "E22" "  -10.13%  "
Set-TextValue "D23" "230.14"
Set-DirectValue "E23" "  -3.63%  "
Set-TextValue "D24" "2.65"
Set-DirectValue "E24" "  -12.01%  "
Set-DirectValue "E25" "  -0.04%  "
Set-DirectValue "E26" "  -6.16%  "
Set-TextValue "D27" "23.44"
Set-DirectValue "E27" "  -5.93%  "
Set-DirectValue "E28" "  -1.24%  "
Set-TextValue "D29" "9.41"
Set-DirectValue "E29" "  -4.20%  "
Set-TextValue "D30" "34.83"
Set-DirectValue "E30" "  -5.52%  "
Set-TextValue "D31" "150.91"
Set-DirectValue "E31" "  -2.72%  "
Set-TextValue "D32" "5.21"
Set-DirectValue "E32" "  -8.49%  "
Set-DirectValue "E33" "  -4.05%  "
Set-DirectValue "E34" "  -3.04%  "
Set-DirectValue "E35" "  -3.59%  "
Set-TextValue "D36" "17.28"
Set-DirectValue "E36" "  -0.64%  "
Set-DirectValue "E37" "  -5.90%  "
Set-DirectValue "E38" "  -7.08%  "
Set-DirectValue "E39" "  -4.02%  "
Set-TextValue "D40" "0.0982"
Set-DirectValue "E40" "  -8.71%  "
Set-TextValue "D41" "3.96"
Set-DirectValue "E41" "  -3.02%  "
Set-DirectValue "E42" "  +1.39%  "
Set-TextValue "D43" "18.87"
Set-DirectValue "E43" "  -11.18%  "
Set-DirectValue "D44" "1.914.80"
Set-DirectValue "E44" "  -4.73%  "
Set-DirectValue "E45" "  -5.59%  "
Set-DirectValue "E46" "  -8.75%  "
Set-TextValue "D47" "8.53"
Set-DirectValue "E47" "  -2.57%  "
Set-DirectValue "D48" "2.674.43"
Set-DirectValue "E48" "  -1.46%  "
Set-TextValue "D49" "93.20"
Set-DirectValue "E49" "  -5.75%  "
Set-DirectValue "E50" "  -7.08%  "
Set-DirectValue "E51" "  -7.92%  "

Write-Host "Applied cryptos update."
